$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the mileage header (M1): mili_meter -> milli_meter
$ws.Range("M1").Value = "milli_meter_reading_on_rent"

# New columns for per-hour / per-km rental pricing
$ws.Range("P1").Value = "per_hour"
$ws.Range("Q1").Value = "per_km"
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0
}

# New car entry row (id 53)
$row = 54
$ws.Cells.Item($row, 1).Value = 53
$ws.Cells.Item($row, 2).Value = "new"
$ws.Cells.Item($row, 3).Value = "no"
$ws.Cells.Item($row, 4).Value = "no"
$ws.Cells.Item($row, 5).Value = "yes"
$ws.Cells.Item($row, 6).Value = 200
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "yes"
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 200
$ws.Cells.Item($row, 17).Value = 200

# Cosmetic window / view refresh
$excel.CalculateFullRebuild()
